$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "1.00" or "63.421.72" must remain text, not be
# auto-converted to numbers/dates by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.421.72'
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D3").Value = '2.591.87'
$ws.Range("E3").Value = '  -2.73%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '572.49'
$ws.Range("E5").Value = '  -4.34%  '
$ws.Range("D6").Value = '154.73'
$ws.Range("E6").Value = '  -2.53%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -4.58%  '
$ws.Range("E9").Value = '  -7.67%  '
$ws.Range("D10").Value = '5.85'
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").Value = '0.379'
$ws.Range("E11").Value = '  -5.83%  '
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").Value = '28.04'
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("D14").Value = '3.059.82'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '0.0000177'
$ws.Range("E15").Value = '  -9.30%  '
$ws.Range("D16").Value = '63.284.33'
$ws.Range("E16").Value = '  -3.71%  '
$ws.Range("D17").Value = '2.587.95'
$ws.Range("E17").Value = '  -2.60%  '
$ws.Range("D18").Value = '11.95'
$ws.Range("E18").Value = '  -5.40%  '
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  -5.78%  '
$ws.Range("D20").Value = '7.47'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '341.27'
$ws.Range("E21").Value = '  -3.10%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '67.33'
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("D24").Value = '1.80'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("E25").Value = '  -4.65%  '
$ws.Range("D26").Value = '9.13'
$ws.Range("E26").Value = '  -5.48%  '
$ws.Range("D27").Value = '577.47'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '0.160'
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("E33").Value = '  -5.23%  '
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("D35").Value = '5.41'
$ws.Range("E35").Value = '  -2.59%  '
$ws.Range("D36").Value = '0.401'
$ws.Range("E36").Value = '  -5.26%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  -4.49%  '
$ws.Range("D39").Value = '154.66'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("E40").Value = '  -5.36%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '2.48'
$ws.Range("E42").Value = '  +6.73%  '
$ws.Range("D43").Value = '41.24'
$ws.Range("E43").Value = '  -3.88%  '
$ws.Range("D44").Value = '156.54'
$ws.Range("E45").Value = '  -5.35%  '
$ws.Range("D46").Value = '23.19'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = '0.0586'
$ws.Range("E47").Value = '  -5.29%  '
$ws.Range("D48").Value = '0.624'
$ws.Range("E48").Value = '  -3.28%  '
$ws.Range("D49").Value = '0.1000'
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("E50").Value = '  -5.17%  '
$ws.Range("D51").Value = '18.74'
$ws.Range("E51").Value = '  -5.44%  '
